$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 151, shifting existing rows 151:249 down to 152:250
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record's data
$ws.Cells.Item(151, 1).Value = 10
$ws.Cells.Item(151, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(151, 3).Value = "La Araucanía"
$ws.Cells.Item(151, 4).Value = 44907
$ws.Cells.Item(151, 5).Value = 9
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100103
$ws.Cells.Item(151, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(151, 9).Value = 100103002
$ws.Cells.Item(151, 10).Value = "Ciruela"
$ws.Cells.Item(151, 11).Value = "Red Beaut"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 230
$ws.Cells.Item(151, 14).Value = 14000
$ws.Cells.Item(151, 15).Value = 14400
$ws.Cells.Item(151, 16).Value = 14165
$ws.Cells.Item(151, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(151, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(151, 19).Value = 787
$ws.Cells.Item(151, 20).Value = 18

# Apply the same number format/style as the other Fecha (date) cells in column D
$ws.Cells.Item(151, 4).NumberFormat = $ws.Cells.Item(152, 4).NumberFormat
